# Regenerate the "K" column (col G, strikeouts) of the save-data sheet for
# gott_trevor with newly computed values (s_vals), per the commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,70,71,72,73)
$s_vals = @(1,1,1,0,1,0,0,2,0,1,1,1,3,2,1,1,0,2,1,1,1,1,2,1,1,0,1,1,1,0,0,1,1,0,1,0,1,2,0,2,2,1,1,0,1,2,2,1,2,2,1,2,0,1,1,0,2,0,2,2,0,1,0,2,1,3,1,1,1,1)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 7).Value = $s_vals[$i]
}
